# Creating-pairs-list-page.html sheet creation + arenas list update on competition page
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the last sheet (group-page.html) to seed the new translation
#    sheet with the same look & feel, then rename/relocate it.
# ---------------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item(6)
$srcSheet.Copy($null, $srcSheet)
$newSheet = $wb.Worksheets.Item(7)
$newSheet.Name = "creating-pairs-list-page.html"

# ---------------------------------------------------------------------------
# 2. Wipe the copied content and rebuild the translation table for the new
#    "creating pairs list" page.
# ---------------------------------------------------------------------------
$newSheet.Cells.Clear()

$rows = @(
    @("origin", "en", "ua", "ru"),
    @("PHX_PAGE_NAME", "Create arena", "Створити арену", "Создать арену"),
    @("PHX_PLACE_NAME", "Arena name", "Назва арени", "Название арены"),
    @("PHX_AUTOSORT_TAB", "Automatic sorting", "Автоматичне сортування", "Автоматическая сортировка"),
    @("PHX_MANUAL_TAB", "Manual sorting", "Ручне сортування", "Ручная сортировка"),
    @("PHX_ACTIVE_GROUPS_LIST", "Selected groups", "Обрані групи", "Выбранные группы"),
    @("PHX_UNACTIVE_GROUPS_LIST", "Groups list", "Список груп", "Список групп"),
    @("PHX_DISTANCE", "One member pairs distance", "Відстань між парами одного учасника", "Расстояние между парами одного участника"),
    @("PHX_AGE_MIN", "Minimal age", "Мінімальний вік", "Минимальный возраст"),
    @("PHX_AGE_MAX", "Maximal age", "Максимальний вік", "Максимальный возраст"),
    @("PHX_WEIGHT_MIN", "Minimal weight", "Мінімальна вага", "Минимальный вес"),
    @("PHX_WEIGHT_MAX", "Maximal weight", "Максимальна вага", "Максимальный вес"),
    @("PHX_FINAL_MIN", "Minimal final part", "Мінімальна фінальна частина", "Минимальная финальная часть"),
    @("PHX_FINAL_MAX", "Maximal final part", "Максимальна фінальна частина", "Максимальная финальная часть"),
    @("PHX_APPLY_BTN", "CREATE", "СТВОРИТИ", "СОЗДАТЬ"),
    @("PHX_ACTIVE_PAIRS_LIST", "Selected pairs", "Обрані пари", "Выбранные пары"),
    @("PHX_UNACTIVE_PAIRS_LIST", "Unattached pairs", "Неприв’язані пари", "Непривязанные пары")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $newSheet.Cells.Item($r + 1, $c + 1).Value = $rowData[$c]
    }
}

$newSheet.Activate()
